$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet had two quantity columns ("原库位数量" in I, "现库位数量" in J)
# inserted between "现库位" (H) and "备注" (old K). The customer asked for
# these to be removed, so delete columns I:J which also shifts the old
# "备注" column (K) left into I, carrying its data/styles with it.
# ---------------------------------------------------------------------------
$ws.Columns("I:J").Delete()

# Update the active selection/cursor.
$ws.Range("M9").Select()

# ---------------------------------------------------------------------------
# Fix up conditional formatting ranges that referenced the removed columns.
# ---------------------------------------------------------------------------

# Rule on row 6 used to span B6:J6 -> now only B6:H6 (J column is gone).
$ruleRow6 = $ws.Range("B6:H6").FormatConditions.Item(1)
$ruleRow6.ModifyAppliesToRange($ws.Range("B6:H6"))

# The big "zebra-stripe" rule used to cover (among others) the old I/J/K
# columns. Reuse the existing rule for the first area (keeps its original
# dxf/priority), then extend the same look (formula + fill color) onto the
# rest of the logical cells using the new column letters (old K -> new I).
$oldBigRule = $ws.Range("A3:G3").FormatConditions.Item(1)
$oldBigRule.ModifyAppliesToRange($ws.Range("A3:G3"))

$bigRuleExtraAreas = @(
    "E4:G4",
    "A4:C4",
    "A5:G5",
    "D2:E5",
    "A8:G17",
    "A6:A7",
    "A2:H2",
    "I2:I17"
)

foreach ($area in $bigRuleExtraAreas) {
    $target = $ws.Range($area)
    $newRule = $target.FormatConditions.Add(2, 3, "=MOD(ROW(),2)=0")
    $newRule.Interior.Color = 16379876
}
